$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 86.360967
$ws.Range("H2").Value = 259.082901
$ws.Range("I2").Value = 0.2482072712525276
$ws.Range("J2").Value = 0.2482072712525276
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 14518.1272507097
$ws.Range("R2").Value = 130663.1452563873
$ws.Range("S2").Value = 0.07406957351864578
$ws.Range("T2").Value = 0.07406957351864581
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 86.360967
$ws.Range("H3").Value = 259.082901
$ws.Range("I3").Value = 0.2482072712525276
$ws.Range("J3").Value = 0.2482072712525276
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 14077.37613920322
$ws.Range("R3").Value = 126696.385252829
$ws.Range("S3").Value = 0.07182091938485882
$ws.Range("T3").Value = 0.07182091938485882
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 86.360967
$ws.Range("H4").Value = 259.082901
$ws.Range("I4").Value = 0.2482072712525276
$ws.Range("J4").Value = 0.2482072712525276
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 14335.36176654351
$ws.Range("R4").Value = 129018.2558988916
$ws.Range("S4").Value = 0.07313712808457945
$ws.Range("T4").Value = 0.07313712808457946
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 86.360967
$ws.Range("H5").Value = 259.082901
$ws.Range("I5").Value = 0.2482072712525276
$ws.Range("J5").Value = 0.2482072712525276
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 5719.404818278767
$ws.Range("R5").Value = 51474.64336450891
$ws.Range("S5").Value = 0.02917965026444354
$ws.Range("T5").Value = 0.02917965026444354
# Row 6
$ws.Range("I6").Value = 0.6003523616657895
$ws.Range("J6").Value = 0.6003523616657896
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 35115.78020234677
$ws.Range("R6").Value = 316042.021821121
$ws.Range("S6").Value = 0.1791560866250972
$ws.Range("T6").Value = 0.1791560866250972
# Row 7
$ws.Range("I7").Value = 0.6003523616657895
$ws.Range("J7").Value = 0.6003523616657896
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.173717145159055
$ws.Range("T7").Value = 0.173717145159055
# Row 8
$ws.Range("I8").Value = 0.6003523616657895
$ws.Range("J8").Value = 0.6003523616657896
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 34673.71543326703
$ws.Range("R8").Value = 312063.4388994033
$ws.Range("S8").Value = 0.1769007303833509
$ws.Range("T8").Value = 0.1769007303833509
# Row 9
$ws.Range("I9").Value = 0.6003523616657895
$ws.Range("J9").Value = 0.6003523616657896
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 13833.83400755786
$ws.Range("R9").Value = 124504.5060680207
$ws.Range("S9").Value = 0.07057839949828651
$ws.Range("T9").Value = 0.07057839949828652
# Row 10
$ws.Range("G10").Value = 52.26262533333333
$ws.Range("H10").Value = 156.787876
$ws.Range("I10").Value = 0.1502063266901572
$ws.Range("J10").Value = 0.1502063266901572
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 8785.860920773357
$ws.Range("R10").Value = 79072.74828696022
$ws.Range("S10").Value = 0.04482430551530035
$ws.Range("T10").Value = 0.04482430551530037
# Row 11
$ws.Range("G11").Value = 52.26262533333333
$ws.Range("H11").Value = 156.787876
$ws.Range("I11").Value = 0.1502063266901572
$ws.Range("J11").Value = 0.1502063266901572
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 8519.133821644036
$ws.Range("R11").Value = 76672.20439479632
$ws.Range("S11").Value = 0.04346349897756949
$ws.Range("T11").Value = 0.0434634989775695
# Row 12
$ws.Range("G12").Value = 52.26262533333333
$ws.Range("H12").Value = 156.787876
$ws.Range("I12").Value = 0.1502063266901572
$ws.Range("J12").Value = 0.1502063266901572
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 8675.257666147425
$ws.Range("R12").Value = 78077.31899532683
$ws.Range("S12").Value = 0.04426002227418767
$ws.Range("T12").Value = 0.04426002227418768
# Row 13
$ws.Range("G13").Value = 52.26262533333333
$ws.Range("H13").Value = 156.787876
$ws.Range("I13").Value = 0.1502063266901572
$ws.Range("J13").Value = 0.1502063266901572
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 3461.183003513203
$ws.Range("R13").Value = 31150.64703161882
$ws.Range("S13").Value = 0.01765849992309967
$ws.Range("T13").Value = 0.01765849992309967
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4293706666666666
$ws.Range("H14").Value = 1.288112
$ws.Range("I14").Value = 0.001234040391525629
$ws.Range("J14").Value = 0.001234040391525629
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 72.18142863533154
$ws.Range("R14").Value = 649.6328577179839
$ws.Range("S14").Value = 0.0003682601442086286
$ws.Range("T14").Value = 0.0003682601442086287
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4293706666666666
$ws.Range("H15").Value = 1.288112
$ws.Range("I15").Value = 0.001234040391525629
$ws.Range("J15").Value = 0.001234040391525629
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 69.99009607902043
$ws.Range("R15").Value = 629.910864711184
$ws.Range("S15").Value = 0.0003570802540560917
$ws.Range("T15").Value = 0.0003570802540560918
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4293706666666666
$ws.Range("H16").Value = 1.288112
$ws.Range("I16").Value = 0.001234040391525629
$ws.Range("J16").Value = 0.001234040391525629
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 71.27275263845333
$ws.Range("R16").Value = 641.4547737460799
$ws.Range("S16").Value = 0.0003636241989249758
$ws.Range("T16").Value = 0.0003636241989249758
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4293706666666666
$ws.Range("H17").Value = 1.288112
$ws.Range("I17").Value = 0.001234040391525629
$ws.Range("J17").Value = 0.001234040391525629
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 28.43581707185956
$ws.Range("R17").Value = 255.922353646736
$ws.Range("S17").Value = 0.000145075794335933
$ws.Range("T17").Value = 0.000145075794335933
